$wb = $excel.ActiveWorkbook

# Update "想去人数" (people interested) counts for two rows on both the
# "展览" sheet and the "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 9334
    $ws.Range("F4").Value = 20
}
